$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a pure percentage (e.g. "79%") need the
# cell pre-formatted as Text; otherwise Excel auto-converts the literal into
# a numeric percentage (0.79) instead of keeping the original text string.
$percentCells = @("H6", "H12", "H24", "H37", "H39", "H46")
foreach ($addr in $percentCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-13 22:18:24"
$ws.Range("I2").Value = "4.3 mm"
$ws.Range("E3").Value = "2026-02-13 22:18:27"
$ws.Range("I3").Value = "7.6 mm"
$ws.Range("E4").Value = "2026-02-13 22:18:29"
$ws.Range("J4").Value = "993.5 hPa"
$ws.Range("L4").Value = "25.2 km/h - 285º 21:49 TU"
$ws.Range("E5").Value = "2026-02-13 22:18:32"
$ws.Range("I5").Value = "3.4 mm"
$ws.Range("N5").Value = "-3.9 °C 21:59 TU"
$ws.Range("E6").Value = "2026-02-13 22:18:34"
$ws.Range("H6").Value = "79%"
$ws.Range("J6").Value = "993.5 hPa"
$ws.Range("O6").Value = "9.0 °C"
$ws.Range("E7").Value = "2026-02-13 22:18:37"
$ws.Range("J7").Value = "993.8 hPa"
$ws.Range("E8").Value = "2026-02-13 22:18:39"
$ws.Range("J8").Value = "993.8 hPa"
$ws.Range("O8").Value = "9.0 °C"
$ws.Range("E9").Value = "2026-02-13 22:18:42"
$ws.Range("O9").Value = "9.7 °C"
$ws.Range("E10").Value = "2026-02-13 22:18:44"
$ws.Range("O10").Value = "8.8 °C"
$ws.Range("E11").Value = "2026-02-13 22:18:47"
$ws.Range("E12").Value = "2026-02-13 22:18:49"
$ws.Range("H12").Value = "84%"
$ws.Range("E13").Value = "2026-02-13 22:18:51"
$ws.Range("E14").Value = "2026-02-13 22:18:54"
$ws.Range("L14").Value = "45.4 km/h - 294º 21:55 TU"
$ws.Range("E15").Value = "2026-02-13 22:18:56"
$ws.Range("I15").Value = "5.1 mm"
$ws.Range("O15").Value = "9.7 °C"
$ws.Range("E16").Value = "2026-02-13 22:18:59"
$ws.Range("I16").Value = "14.4 mm"
$ws.Range("E17").Value = "2026-02-13 22:19:01"
$ws.Range("E18").Value = "2026-02-13 22:19:04"
$ws.Range("J18").Value = "993.7 hPa"
$ws.Range("E19").Value = "2026-02-13 22:19:06"
$ws.Range("E20").Value = "2026-02-13 22:19:09"
$ws.Range("I20").Value = "24.3 mm"
$ws.Range("E21").Value = "2026-02-13 22:19:11"
$ws.Range("J21").Value = "996.7 hPa"
$ws.Range("N21").Value = "-0.2 °C 21:52 TU"
$ws.Range("O21").Value = "0.9 °C"
$ws.Range("E22").Value = "2026-02-13 22:19:14"
$ws.Range("L22").Value = "56.9 km/h - 312º 21:57 TU"
$ws.Range("E23").Value = "2026-02-13 22:19:16"
$ws.Range("G23").Value = "190 cm"
$ws.Range("I23").Value = "13.2 mm"
$ws.Range("L23").Value = "78.8 km/h - 272º 21:34 TU"
$ws.Range("E24").Value = "2026-02-13 22:19:19"
$ws.Range("H24").Value = "94%"
$ws.Range("J24").Value = "994.8 hPa"
$ws.Range("L24").Value = "65.5 km/h - 296º 21:50 TU"
$ws.Range("E25").Value = "2026-02-13 22:19:21"
$ws.Range("I25").Value = "9.8 mm"
$ws.Range("E26").Value = "2026-02-13 22:19:24"
$ws.Range("E27").Value = "2026-02-13 22:19:26"
$ws.Range("E28").Value = "2026-02-13 22:19:29"
$ws.Range("J28").Value = "993.9 hPa"
$ws.Range("E29").Value = "2026-02-13 22:19:31"
$ws.Range("E30").Value = "2026-02-13 22:19:34"
$ws.Range("J30").Value = "993.4 hPa"
$ws.Range("E31").Value = "2026-02-13 22:19:36"
$ws.Range("I31").Value = "4.9 mm"
$ws.Range("J31").Value = "992.4 hPa"
$ws.Range("O31").Value = "10.1 °C"
$ws.Range("E32").Value = "2026-02-13 22:19:39"
$ws.Range("E33").Value = "2026-02-13 22:19:41"
$ws.Range("J33").Value = "995.5 hPa"
$ws.Range("E34").Value = "2026-02-13 22:19:44"
$ws.Range("E35").Value = "2026-02-13 22:19:47"
$ws.Range("I35").Value = "8.8 mm"
$ws.Range("O35").Value = "5.8 °C"
$ws.Range("E36").Value = "2026-02-13 22:19:49"
$ws.Range("J36").Value = "993.5 hPa"
$ws.Range("E37").Value = "2026-02-13 22:19:52"
$ws.Range("H37").Value = "86%"
$ws.Range("J37").Value = "995.4 hPa"
$ws.Range("E38").Value = "2026-02-13 22:19:54"
$ws.Range("N38").Value = "7.9 °C 21:46 TU"
$ws.Range("O38").Value = "9.6 °C"
$ws.Range("E39").Value = "2026-02-13 22:19:56"
$ws.Range("H39").Value = "80%"
$ws.Range("E40").Value = "2026-02-13 22:19:59"
$ws.Range("J40").Value = "997.1 hPa"
$ws.Range("O40").Value = "1.6 °C"
$ws.Range("E41").Value = "2026-02-13 22:20:01"
$ws.Range("J41").Value = "994.2 hPa"
$ws.Range("E42").Value = "2026-02-13 22:20:04"
$ws.Range("E43").Value = "2026-02-13 22:20:06"
$ws.Range("E44").Value = "2026-02-13 22:20:08"
$ws.Range("I44").Value = "10.6 mm"
$ws.Range("E45").Value = "2026-02-13 22:20:11"
$ws.Range("I45").Value = "2.6 mm"
$ws.Range("E46").Value = "2026-02-13 22:20:13"
$ws.Range("H46").Value = "87%"
$ws.Range("J46").Value = "995.0 hPa"
$ws.Range("L46").Value = "67.7 km/h - 303º 21:32 TU"
$ws.Range("O46").Value = "9.2 °C"
